$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recording Sheet")

# Remove the "Test Step Number" / "Run" columns (old D:E) and the
# "Test Step Detail" column (old G, now E after the first delete).
$ws.Range("D:E").Delete()
$ws.Range("E:E").Delete()

# Re-apply the autofilter over the new, narrower header range (A1:L1).
$ws.AutoFilterMode = $false
$ws.Range("A1:L1").AutoFilter()

# Populate the new defect-log rows (2-9).
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "UAT_analysis.Rmd"
$ws.Range("C2").Value = "Transformation of data from FTP S3 bucket to REL input bucket"
$ws.Range("E2").Value = "Original multi-SLX data not copied from 'ftps-core-ftps3bucket-1b6k14wlq0yu7' to 'transformed-data-bucket-rel' (e.g.  'SLX-100.s_1.SRR8983578.fastq.gz'). Possibly because the filenames are not as expected - need clarification."
$ws.Range("F2").Value = 43872
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").Value = "ATK"
$ws.Range("H2").Value = "Major"
$ws.Range("I2").Value = "Identified"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "UAT_analysis.Rmd"
$ws.Range("C3").Value = "Testing concordance of v1 and v2 output"
$ws.Range("E3").Value = "The order of the columns (sample names) in the 'combined_counts.txt' file seems to depend on the order of samples in the original json analysis config file even though the data is identical."
$ws.Range("F3").Value = 43872
$ws.Range("F3").NumberFormat = "mm-dd-yy"
$ws.Range("G3").Value = "ATK"
$ws.Range("H3").Value = "Trivial"
$ws.Range("I3").Value = "Closed"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "UAT_analysis.Rmd"
$ws.Range("C4").Value = "Testing concordance of v1 and v2 output"
$ws.Range("E4").Value = "The order of the sgRNAs in the 'combined_counts.txt' file also seems to depend on the order of the samples in the original json file (perhaps ordering of rows occurs across sample column order). "
$ws.Range("F4").Value = 43872
$ws.Range("F4").NumberFormat = "mm-dd-yy"
$ws.Range("G4").Value = "ATK"
$ws.Range("H4").Value = "Trivial"
$ws.Range("I4").Value = "Closed"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "UAT_analysis.Rmd"
$ws.Range("C5").Value = "Testing concordance of v1 and v2 output"
$ws.Range("E5").Value = "The 'neg.rank' and 'pos.rank' columns are not identical in the Mageck output - there seems to be some randomness in the ranking of tied genes. This also impacts the order of the genes in the output file."
$ws.Range("F5").Value = 43872
$ws.Range("F5").NumberFormat = "mm-dd-yy"
$ws.Range("G5").Value = "ATK"
$ws.Range("H5").Value = "Trivial"
$ws.Range("I5").Value = "Closed"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "UAT_analysis.Rmd"
$ws.Range("C6").Value = "Testing concordance of v1 and v2 output"
$ws.Range("E6").Value = "'" + '"GeneConnectivity" interpretation output is private to AZ'
$ws.Range("F6").Value = 43872
$ws.Range("F6").NumberFormat = "mm-dd-yy"
$ws.Range("G6").Value = "ATK"
$ws.Range("H6").Value = "Trivial"
$ws.Range("I6").Value = "Closed"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "UAT_analysis.Rmd"
$ws.Range("C7").Value = "Testing concordance of v1 and v2 output"
$ws.Range("E7").Value = "'" + '"PathwayAnalysis" interpretation output is provate to v2.'
$ws.Range("F7").Value = 43872
$ws.Range("F7").NumberFormat = "mm-dd-yy"
$ws.Range("G7").Value = "ATK"
$ws.Range("H7").Value = "Trivial"
$ws.Range("I7").Value = "Closed"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "UAT_analysis.Rmd"
$ws.Range("C8").Value = "Testing concordance of v1 and v2 output"
$ws.Range("E8").Value = "`"fastq`" and `"fastqc`" sub-folders no longer saved in v2 'crisprn' sub-folder, but at the same level as 'crisprn' folder."
$ws.Range("F8").Value = 43872
$ws.Range("F8").NumberFormat = "mm-dd-yy"
$ws.Range("G8").Value = "ATK"
$ws.Range("H8").Value = "Trivial"
$ws.Range("I8").Value = "Closed"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "UAT_analysis.Rmd"
$ws.Range("C9").Value = "Testing concordance of v1 and v2 output"
$ws.Range("E9").Value = "`"counts' sub-folder within `"crisprn`" folder now contains the merging counts log, for v1 it contained the log for individual sample counts as well as the merging counts log."
$ws.Range("F9").Value = 43872
$ws.Range("F9").NumberFormat = "mm-dd-yy"
$ws.Range("G9").Value = "ATK"
$ws.Range("H9").Value = "Trivial"
$ws.Range("I9").Value = "Closed"

# Wrap text for the whole data block so long descriptions match the target formatting.
$ws.Range("A2:I9").WrapText = $true

# Row heights: first data row taller (wraps a long line twice), remaining
# rows share a slightly shorter height.
$ws.Rows.Item(2).RowHeight = 80
$ws.Range("3:9").RowHeight = 64
